$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K" header) values per regen of save_data
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
